$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Simple single occurrence line-number updates in the stack trace
Replace-Text "JavaMethodService.java:163" "JavaMethodService.java:162"
Replace-Text "AbstractService.java:136" "AbstractService.java:135"
Replace-Text "EvaluationServices.java:168" "EvaluationServices.java:172"
Replace-Text "EvaluationServices.java:204" "EvaluationServices.java:208"
Replace-Text "AstEvaluator.java:192" "AstEvaluator.java:189"
Replace-Text "AstSwitch.java:118" "AstSwitch.java:119"
Replace-Text "AstEvaluator.java:112" "AstEvaluator.java:109"
Replace-Text "GeneratedMethodAccessor74" "GeneratedMethodAccessor73"

# Replace the tail of the stack trace (from the surefire/equinox launcher chain)
# with the new Eclipse JDT JUnit runner chain.
$oldTail = "	at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)" + `
"`n	at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)" + `
"`n	at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)" + `
"`n	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + `
"`n	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + `
"`n	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + `
"`n	at java.lang.reflect.Method.invoke(Method.java:498)" + `
"`n	at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)" + `
"`n	at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)" + `
"`n	at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)" + `
"`n	at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)" + `
"`n	at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)" + `
"`n	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + `
"`n	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + `
"`n	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + `
"`n	at java.lang.reflect.Method.invoke(Method.java:498)" + `
"`n	at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)" + `
"`n	at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)" + `
"`n	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)" + `
"`n	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)" + `
"`n	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)" + `
"`n	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)" + `
"`n	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + `
"`n	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + `
"`n	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + `
"`n	at java.lang.reflect.Method.invoke(Method.java:498)" + `
"`n	at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)" + `
"`n	at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)" + `
"`n	at org.eclipse.equinox.launcher.Main.run(Main.java:1498)" + `
"`n	at org.eclipse.equinox.launcher.Main.main(Main.java:1471)"

$newTail = "	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)" + `
"`n	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)" + `
"`n	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)" + `
"`n	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)" + `
"`n	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)" + `
"`n	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

Replace-Text $oldTail $newTail
